$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that used to hold only the "5840938 - Marcelo Rodrigues de Holanda"
# text in columns B/C (with no label in column A) is removed; everything
# below it shifts up by one row.
$ws.Rows.Item(13).Delete()

# Column A keeps its labels (now shifted), but several B/C "answer" cells
# end up with different text than what a plain row-shift would produce.
$ws.Range("B10").Value = "5840938 - Marcelo Rodrigues de Holanda"
$ws.Range("C10").Value = "5840938 - Marcelo Rodrigues de Holanda"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2012" looks like a date, so a plain assignment would get
# auto-converted to a date serial number. Force it in as text (matching
# the original shared-string cell) and then restore the normal
# column B/C number format via a format-only paste from a sibling cell.
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("B10").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "01/01/2012"
$ws.Range("C10").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("B18").Value = "5840938 - Marcelo Rodrigues de Holanda"
$ws.Range("C18").Value = "5840938 - Marcelo Rodrigues de Holanda"

$ws.Range("B19").Value = "Aula expositiva e exercícios dirigidos."
$ws.Range("C19").Value = "Aula expositiva e exercícios dirigidos."

$ws.Range("B20").Value = "Média ponderada de exercícios e provas."
$ws.Range("C20").Value = "Média ponderada de exercícios e provas."

$ws.Range("B21").Value = "Prova única com nota igual ou superior a 5,0."
$ws.Range("C21").Value = "Prova única com nota igual ou superior a 5,0."
